$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 956.0185
$ws.Range("I15").Value = 956.0185
$ws.Range("K15").Value = 2868.0555
$ws.Range("M15").Value = -2699.0555
$ws.Range("H19").Value = 2193961.8
$ws.Range("I19").Value = 4386661.5
$ws.Range("J19").Value = 1261.8334
$ws.Range("K19").Value = 4386661.5
$ws.Range("L19").Value = 1261.8334
$ws.Range("M19").Value = -4386486.5
$ws.Range("N19").Value = -1611.8334
$ws.Range("H34").Value = 4036.3635
$ws.Range("I34").Value = 4036.3635
$ws.Range("K34").Value = 4036.3635
$ws.Range("M34").Value = -3833.3635
$ws.Range("H36").Value = 4036.3635
$ws.Range("I36").Value = 4036.3635
$ws.Range("K36").Value = 4036.3635
$ws.Range("M36").Value = -3321.3635
$ws.Range("H116").Value = 633975.0600000001
$ws.Range("J116").Value = 12130.6
$ws.Range("L116").Value = 12130.6
$ws.Range("N116").Value = -19014.6
$ws.Range("H121").Value = 2955
$ws.Range("J121").Value = 2955
$ws.Range("L121").Value = 8865
$ws.Range("N121").Value = -12359
$ws.Range("H132").Value = 33339048
$ws.Range("I132").Value = 43484108
$ws.Range("J132").Value = 5285.7144
$ws.Range("K132").Value = 130452324
$ws.Range("L132").Value = 15857.1432
$ws.Range("M132").Value = -130449794
$ws.Range("N132").Value = -20917.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1527.1428
$ws.Range("I61").Value = 938
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 938
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -726
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 6118.579
$ws.Range("I74").Value = 7004.5386
$ws.Range("J74").Value = 4199
$ws.Range("K74").Value = 7004.5386
$ws.Range("L74").Value = 4199
$ws.Range("M74").Value = -6130.5386
$ws.Range("N74").Value = -5947
$ws.Range("H77").Value = 6118.579
$ws.Range("I77").Value = 7004.5386
$ws.Range("J77").Value = 4199
$ws.Range("K77").Value = 35022.693
$ws.Range("L77").Value = 20995
$ws.Range("M77").Value = -30654.693
$ws.Range("N77").Value = -29731
$ws.Range("H132").Value = 2190.16
$ws.Range("I132").Value = 1044.875
$ws.Range("J132").Value = 4226.222
$ws.Range("K132").Value = 3134.625
$ws.Range("L132").Value = 12678.666
$ws.Range("M132").Value = -604.625
$ws.Range("N132").Value = -17738.666
$ws.Range("H133").Value = 22425.666
$ws.Range("J133").Value = 22425.666
$ws.Range("L133").Value = 22425.666
$ws.Range("N133").Value = -27485.666
$ws.Range("H136").Value = 1527.1428
$ws.Range("I136").Value = 938
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2814
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -264
$ws.Range("N136").Value = -14100
$ws.Range("H137").Value = 41139.668
$ws.Range("J137").Value = 41139.668
$ws.Range("L137").Value = 41139.668
$ws.Range("N137").Value = -51339.668
$ws.Range("H139").Value = 41206
$ws.Range("J139").Value = 41206
$ws.Range("L139").Value = 41206
$ws.Range("N139").Value = -51486

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 45750
$ws.Range("J137").Value = 45750
$ws.Range("L137").Value = 45750
$ws.Range("N137").Value = -55950

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 323.875
$ws.Range("I7").Value = 403
$ws.Range("J7").Value = 284.3125
$ws.Range("K7").Value = 403
$ws.Range("L7").Value = 284.3125
$ws.Range("M7").Value = -290
$ws.Range("N7").Value = -510.3125
$ws.Range("H99").Value = 7146139.5
$ws.Range("I99").Value = 16668273
$ws.Range("J99").Value = 4538.75
$ws.Range("K99").Value = 16668273
$ws.Range("L99").Value = 4538.75
$ws.Range("M99").Value = -16666775
$ws.Range("N99").Value = -7534.75
$ws.Range("H115").Value = 25000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
$ws.Range("H126").Value = 7146139.5
$ws.Range("I126").Value = 16668273
$ws.Range("J126").Value = 4538.75
$ws.Range("K126").Value = 50004819
$ws.Range("L126").Value = 13616.25
$ws.Range("M126").Value = -50002349
$ws.Range("N126").Value = -18556.25
$ws.Range("H132").Value = 2933.8333
$ws.Range("I132").Value = 1789.125
$ws.Range("J132").Value = 5223.25
$ws.Range("K132").Value = 5367.375
$ws.Range("L132").Value = 15669.75
$ws.Range("M132").Value = -2837.375
$ws.Range("N132").Value = -20729.75
$ws.Range("H134").Value = 5134.1724
$ws.Range("I134").Value = 5909.65
$ws.Range("K134").Value = 17728.95
$ws.Range("M134").Value = -15193.95
$ws.Range("H137").Value = 45590
$ws.Range("J137").Value = 45590
$ws.Range("L137").Value = 45590
$ws.Range("N137").Value = -55790
$ws.Range("H138").Value = 44538.57
$ws.Range("J138").Value = 44538.57
$ws.Range("L138").Value = 44538.57
$ws.Range("N138").Value = -54818.57
$ws.Range("H140").Value = 87930
$ws.Range("J140").Value = 87930
$ws.Range("L140").Value = 87930
$ws.Range("N140").Value = -98290
$ws.Range("H141").Value = 25271.428
$ws.Range("J141").Value = 25271.428
$ws.Range("L141").Value = 25271.428
$ws.Range("N141").Value = -35631.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 446380.72
$ws.Range("I5").Value = 701.9
$ws.Range("J5").Value = 669220.1
$ws.Range("K5").Value = 2105.7
$ws.Range("L5").Value = 2007660.3
$ws.Range("M5").Value = -1993.7
$ws.Range("N5").Value = -2007884.3
$ws.Range("H82").Value = 1668.8334
$ws.Range("I82").Value = 1003.25
$ws.Range("K82").Value = 3009.75
$ws.Range("M82").Value = -2603.75
$ws.Range("H85").Value = 1668.8334
$ws.Range("I85").Value = 1003.25
$ws.Range("K85").Value = 3009.75
$ws.Range("M85").Value = -1605.75
$ws.Range("H113").Value = 4032857
$ws.Range("I113").Value = 631
$ws.Range("J113").Value = 7813069
$ws.Range("K113").Value = 1893
$ws.Range("L113").Value = 23439207
$ws.Range("M113").Value = 277
$ws.Range("N113").Value = -23443547
$ws.Range("H121").Value = 1039.375
$ws.Range("I121").Value = 538
$ws.Range("J121").Value = 1069.5784
$ws.Range("K121").Value = 1614
$ws.Range("L121").Value = 3208.7352
$ws.Range("M121").Value = -304
$ws.Range("N121").Value = -5828.7352
$ws.Range("H122").Value = 2959.4614
$ws.Range("J122").Value = 3582.6553
$ws.Range("L122").Value = 32243.8977
$ws.Range("N122").Value = -37143.8977
$ws.Range("H131").Value = 709.09
$ws.Range("I131").Value = 281.35294
$ws.Range("J131").Value = 796.6988
$ws.Range("K131").Value = 844.05882
$ws.Range("L131").Value = 2390.0964
$ws.Range("M131").Value = 4195.94118
$ws.Range("N131").Value = -12470.0964
$ws.Range("H135").Value = 446380.72
$ws.Range("I135").Value = 701.9
$ws.Range("J135").Value = 669220.1
$ws.Range("K135").Value = 6317.099999999999
$ws.Range("L135").Value = 6022980.899999999
$ws.Range("M135").Value = -3782.099999999999
$ws.Range("N135").Value = -6028050.899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 33501.168
$ws.Range("J46").Value = 33501.168
$ws.Range("L46").Value = 33501.168
$ws.Range("N46").Value = -33813.168
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
$ws.Range("H126").Value = 3076.58
$ws.Range("I126").Value = 2877.9873
$ws.Range("J126").Value = 3823.6667
$ws.Range("K126").Value = 8633.961899999998
$ws.Range("L126").Value = 11471.0001
$ws.Range("M126").Value = -6163.961899999998
$ws.Range("N126").Value = -16411.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 4244
$ws.Range("J17").Value = 4618
$ws.Range("L17").Value = 4618
$ws.Range("N17").Value = -4958
$ws.Range("H22").Value = 2087.5833
$ws.Range("J22").Value = 2633.5
$ws.Range("L22").Value = 2633.5
$ws.Range("N22").Value = -3223.5
$ws.Range("H27").Value = 2087.5833
$ws.Range("J27").Value = 2633.5
$ws.Range("L27").Value = 2633.5
$ws.Range("N27").Value = -2847.5
$ws.Range("H106").Value = 39999.5
$ws.Range("J106").Value = 39999.5
$ws.Range("L106").Value = 39999.5
$ws.Range("N106").Value = -42523.5
$ws.Range("H132").Value = 4686.8
$ws.Range("I132").Value = 1667.5
$ws.Range("J132").Value = 6699.6665
$ws.Range("K132").Value = 5002.5
$ws.Range("L132").Value = 20098.9995
$ws.Range("M132").Value = -2472.5
$ws.Range("N132").Value = -25158.9995
$ws.Range("H136").Value = 3589.625
$ws.Range("I136").Value = 1183.4
$ws.Range("J136").Value = 7600
$ws.Range("K136").Value = 3550.2
$ws.Range("L136").Value = 22800
$ws.Range("M136").Value = -1000.2
$ws.Range("N136").Value = -27900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2704324
$ws.Range("I96").Value = 125874.75
$ws.Range("K96").Value = 125874.75
$ws.Range("M96").Value = -124501.75
